# Insert two new rows right after row 49 (shifting existing rows 50-71 down to 52-73)
# and populate them with the new weekly price records (date 2023-08-03 = serial 45141).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(50).Insert()
$ws.Rows.Item(50).Insert()

# New row 50: Membrillo Champion, "Primera" quality
$ws.Range("A50").Value = 7
$ws.Range("B50").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C50").Value = "Ñuble"
$ws.Range("D50").Value = 45141
$ws.Range("E50").Value = 16
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100104
$ws.Range("H50").Value = "Frutos de pepita"
$ws.Range("I50").Value = 100104003
$ws.Range("J50").Value = "Membrillo"
$ws.Range("K50").Value = "Champion"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 40
$ws.Range("N50").Value = 14000
$ws.Range("O50").Value = 14000
$ws.Range("P50").Value = 14000
$ws.Range("Q50").Value = "$/bandeja 18 kilos granel"
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 778
$ws.Range("T50").Value = 18

# New row 51: Membrillo Champion, "Segunda" quality
$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = 45141
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100104
$ws.Range("H51").Value = "Frutos de pepita"
$ws.Range("I51").Value = 100104003
$ws.Range("J51").Value = "Membrillo"
$ws.Range("K51").Value = "Champion"
$ws.Range("L51").Value = "Segunda"
$ws.Range("M51").Value = 40
$ws.Range("N51").Value = 12000
$ws.Range("O51").Value = 12000
$ws.Range("P51").Value = 12000
$ws.Range("Q51").Value = "$/bandeja 18 kilos granel"
$ws.Range("R51").Value = "Región de O'Higgins"
$ws.Range("S51").Value = 667
$ws.Range("T51").Value = 18
